$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.875.78'
$ws.Range('E2').Value = '  -1.02%  '
$ws.Range('D3').Value = '1.638.03'
$ws.Range('E3').Value = '  -0.83%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.42'
$ws.Range('E5').Value = '  -0.11%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5028'
$ws.Range('E6').Value = '  -1.66%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2565'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06388'
$ws.Range('E9').Value = '  -0.57%  '
$ws.Range('E10').Value = '  -1.34%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07725'
$ws.Range('E11').Value = '  -0.79%  '
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.273'
$ws.Range('E12').Value = '  -0.32%  '
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').Value = '1.864.61'
$ws.Range('E13').Value = '  -0.80%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.620.43'
$ws.Range('E14').Value = '  -1.93%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5459'
$ws.Range('E15').Value = '  -1.13%  '
$ws.Range('D16').Value = '0.0₅7914'
$ws.Range('E16').Value = '  -1.37%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '64.05'
$ws.Range('E17').Value = '  -0.17%  '
$ws.Range('D18').Value = '25.901.98'
$ws.Range('E18').Value = '  -0.90%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.003'
$ws.Range('E19').Value = '  -0.13%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '202.71'
$ws.Range('E20').Value = '  -3.58%  '
$ws.Range('E21').Value = '  -0.76%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.923'
$ws.Range('E22').Value = '  -1.44%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.979'
$ws.Range('E23').Value = '  -1.23%  '
$ws.Range('E24').Value = '  +0.00%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.921'
$ws.Range('E25').Value = '  +10.13%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '141.88'
$ws.Range('E26').Value = '  -1.26%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1138'
$ws.Range('E27').Value = '  -3.47%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.65'
$ws.Range('E28').Value = '  -0.98%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.719'
$ws.Range('E29').Value = '  -3.88%  '
$ws.Range('E30').Value = '  -0.07%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.04987'
$ws.Range('E31').Value = '  -2.88%  '
$ws.Range('E32').Value = '  -1.98%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.197'
$ws.Range('E33').Value = '  -0.70%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.374'
$ws.Range('E35').Value = '  +0.68%  '
$ws.Range('D36').Value = '1.166.13'
$ws.Range('E36').Value = '  +0.16%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.626'
$ws.Range('E37').Value = '  -4.38%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.8918'
$ws.Range('E38').Value = '  -3.64%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.5598'
$ws.Range('E39').Value = '  -1.77%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01563'
$ws.Range('E40').Value = '  -1.35%  '
$ws.Range('E41').Value = '  -0.05%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.690'
$ws.Range('E42').Value = '  +0.62%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8067'
$ws.Range('E43').Value = '  -2.78%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '99.72'
$ws.Range('E44').Value = '  -0.37%  '
$ws.Range('D45').Value = '1.776.79'
$ws.Range('E45').Value = '  -0.74%  '
$ws.Range('D46').Value = '0.0₈116'
$ws.Range('E46').Value = '  -0.55%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4523'
$ws.Range('E47').Value = '  -0.59%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.003'
$ws.Range('E48').Value = '  -0.25%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '54.93'
$ws.Range('E49').Value = '  -1.10%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05074'
$ws.Range('E50').Value = '  +0.27%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.005'
$ws.Range('E51').Value = '  -0.14%  '
